# Fix latency units in report sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update header label X2: "Utility" -> "Utility (Percent)"
$ws.Range("X2").Value = "Utility (Percent)"

# Append " usec" unit suffix to Read/Write Latency columns (L:Q) for data rows 3-14
$latencyColumns = @("L", "M", "N", "O", "P", "Q")
for ($row = 3; $row -le 14; $row++) {
    foreach ($col in $latencyColumns) {
        $cell = $ws.Range("$col$row")
        $current = $cell.Value()
        $cell.Value = "$current usec"
    }
}
